$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 276
$ws.Range("C2").Value = 95.5
$ws.Range("B3").Value = 13
$ws.Range("C3").Value = 4.5
